$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "95.739.01"
$ws.Range("E2").Value = "  +2.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.610.06"
$ws.Range("E3").Value = "  +5.37%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.41"
$ws.Range("E5").Value = "  +2.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "655.79"
$ws.Range("E6").Value = "  +5.84%  "
$ws.Range("E7").Value = "  +6.37%  "
$ws.Range("E8").Value = "  +2.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E10").Value = "  +3.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.609.63"
$ws.Range("E11").Value = "  +5.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.06"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.201"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.35"
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.277.83"
$ws.Range("E15").Value = "  +5.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.617.74"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000256"
$ws.Range("E17").Value = "  +3.63%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.24"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.606.44"
$ws.Range("E19").Value = "  +5.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.69"
$ws.Range("E20").Value = "  +9.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.04"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.63"
$ws.Range("E22").Value = "  +7.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.493"
$ws.Range("E23").Value = "  +8.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "510.87"
$ws.Range("E24").Value = "  +1.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000197"
$ws.Range("E25").Value = "  +6.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.67"
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "96.98"
$ws.Range("E27").Value = "  +2.04%  "
$ws.Range("E28").Value = "  +6.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.801.52"
$ws.Range("E29").Value = "  +5.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.12"
$ws.Range("E30").Value = "  +14.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.41"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.178"
$ws.Range("E35").Value = "  +3.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.05"
$ws.Range("E36").Value = "  +7.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.563"
$ws.Range("E37").Value = "  +2.36%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.24"
$ws.Range("E38").Value = "  +10.08%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "575.01"
$ws.Range("E39").Value = "  +3.42%  "
$ws.Range("E40").Value = "  +5.99%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.925"
$ws.Range("E43").Value = "  +1.28%  "
$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.73"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.73"
$ws.Range("E45").Value = "  +4.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.77"
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "34.30"
$ws.Range("E47").Value = "  +35.51%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.25"
$ws.Range("E48").Value = "  +6.44%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0417"
$ws.Range("E49").Value = "  +1.88%  "
$ws.Range("B50").Value = "MantraDAO"
$ws.Range("C50").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.54"
$ws.Range("E50").Value = "  -3.98%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.90"
$ws.Range("E51").Value = "  +0.45%  "
